$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 (was "Good Morning") to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection on the sheet
$ws.Range("E8").Select()
